$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crawl snapshot.
# Values that look like plain numbers are force-written as Text so Excel
# keeps the original author-intended string formatting (e.g. "472.73",
# trailing zeros, multi-dot "thousand.thousand.cents" style prices, etc.)
# rather than silently coercing them into floating-point numbers.

$ws.Range("D2").Value = "54.395.07"
$ws.Range("E2").Value = "  -7.09%  "
$ws.Range("D3").Value = "2.879.73"
$ws.Range("E3").Value = "  -9.66%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "472.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -11.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.92%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "2.878.50"
$ws.Range("E8").Value = "  -9.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.405"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -11.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -10.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0963"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -14.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.330"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -16.24%  "
$ws.Range("E13").Value = "  -4.67%  "
$ws.Range("D14").Value = "3.348.03"
$ws.Range("E14").Value = "  -10.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -9.79%  "
$ws.Range("D16").Value = "54.359.46"
$ws.Range("E16").Value = "  -7.22%  "
$ws.Range("D17").Value = "2.866.04"
$ws.Range("E17").Value = "  -10.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000134"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -14.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -12.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -12.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "297.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -16.79%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.443"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -14.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "58.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -15.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.152"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -10.61%  "
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").Value = "0.0₃0809"
$ws.Range("E29").Value = "  -14.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -12.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -15.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -12.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -15.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "136.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -15.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -14.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -13.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -10.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0621"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.84%  "
$ws.Range("D41").Value = "2.894.48"
$ws.Range("E41").Value = "  -10.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "34.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -14.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.954"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -12.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.602"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -15.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -11.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -15.14%  "
$ws.Range("D48").Value = "2.039.48"
$ws.Range("E48").Value = "  -10.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -14.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -12.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0215"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.82%  "

Write-Host "Updated cryptos list"
